$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as plain text (no numeric/percent/date auto-conversion),
# then strip any style change introduced by the forced quote-prefix so the cell
# keeps its original (default) style, matching the source data format.
function Set-TextValue {
    param($Sheet, $RowNum, $ColNum, $Text)
    $cell = $Sheet.Cells.Item($RowNum, $ColNum)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

# Row-level updates taken from the refreshed symbol-list data source.
# Each entry updates Price (D), Volume 1h (E) when present, and Hora (G).
$updates = @(
    @{Row=2; D='306.91'; E='-3.65%'; G='19'},
    @{Row=3; D=$null; E='-6.14%'; G='19'},
    @{Row=4; D='5.116'; E='-0.41%'; G='19'},
    @{Row=5; D='0.07749'; E='-5.73%'; G='19'},
    @{Row=6; D='1.968'; E='-4.44%'; G='19'},
    @{Row=7; D='4.403'; E='1.87%'; G='19'},
    @{Row=8; D='8.265'; E='-0.71%'; G='19'},
    @{Row=9; D='3.086'; E='-8.25%'; G='19'},
    @{Row=10; D='0.9236'; E='-1.37%'; G='19'},
    @{Row=11; D='0.1308'; E='-2.17%'; G='19'},
    @{Row=12; D='0.1944'; E='-2.38%'; G='19'},
    @{Row=13; D='0.08899'; E='-2.46%'; G='19'},
    @{Row=14; D='0.03449'; E='-2.17%'; G='19'},
    @{Row=15; D='0.09711'; E='-1.04%'; G='19'},
    @{Row=16; D='0.001385'; E='-0.72%'; G='19'},
    @{Row=17; D='0.006151'; E='-2.94%'; G='19'},
    @{Row=18; D='3.592'; E='-2.66%'; G='19'},
    @{Row=19; D=$null; E='-2.04%'; G='19'},
    @{Row=20; D=$null; E='-0.71%'; G='19'},
    @{Row=21; D='5.026'; E='3.56%'; G='19'},
    @{Row=22; D='0.2490'; E='1.67%'; G='19'},
    @{Row=23; D=$null; E='5,172.51%'; G='19'},
    @{Row=24; D='0.04352'; E='0.62%'; G='19'},
    @{Row=25; D='0.001216'; E='-0.87%'; G='19'},
    @{Row=26; D='0.004540'; E='-5.18%'; G='19'},
    @{Row=27; D='0.0001356'; E='4.57%'; G='19'},
    @{Row=28; D=$null; E=$null; G='19'},
    @{Row=29; D=$null; E=$null; G='19'},
    @{Row=30; D=$null; E=$null; G='19'},
    @{Row=31; D=$null; E=$null; G='19'},
    @{Row=32; D=$null; E=$null; G='19'},
    @{Row=33; D=$null; E=$null; G='19'},
    @{Row=34; D=$null; E=$null; G='19'},
    @{Row=35; D=$null; E=$null; G='19'},
    @{Row=36; D=$null; E=$null; G='19'},
    @{Row=37; D=$null; E=$null; G='19'},
    @{Row=38; D=$null; E=$null; G='19'},
    @{Row=39; D='0.02255'; E='1.81%'; G='19'},
    @{Row=40; D='0.04983'; E='-4.42%'; G='19'},
    @{Row=41; D='0.007579'; E='-1.11%'; G='19'},
    @{Row=42; D='0.009837'; E='1.33%'; G='19'},
    @{Row=43; D='0.1350'; E='-3.38%'; G='19'},
    @{Row=44; D='0.002002'; E='-4.13%'; G='19'},
    @{Row=45; D=$null; E='-5.68%'; G='19'},
    @{Row=46; D='0.00006842'; E='2.90%'; G='19'},
    @{Row=47; D=$null; E='0.36%'; G='19'},
    @{Row=48; D='0.003029'; E='5.20%'; G='19'},
    @{Row=49; D='0.001302'; E='-22.82%'; G='19'},
    @{Row=50; D='0.00002103'; E='0.36%'; G='19'},
    @{Row=51; D='0.0002002'; E='0.36%'; G='19'}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) { Set-TextValue $ws $u.Row 4 $u.D }
    if ($null -ne $u.E) { Set-TextValue $ws $u.Row 5 $u.E }
    if ($null -ne $u.G) { Set-TextValue $ws $u.Row 7 $u.G }
}
